$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns A and B, and columns D and E, for every used row (1-7)
for ($r = 1; $r -le 7; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value2 = $b
    $ws.Cells.Item($r, 2).Value2 = $a

    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value2 = $e
    $ws.Cells.Item($r, 5).Value2 = $d
}
